$wb = $excel.ActiveWorkbook

# qol_norm sheet: update column F (US quality-of-life norms)
$qol = $wb.Worksheets.Item("qol_norm")
$qol.Range("F4").Value = 0.88100000000000001
$qol.Range("F5").Value = 0.878
$qol.Range("F6").Value = 0.85499999999999998
$qol.Range("F7").Value = 0.83899999999999997
$qol.Range("F8").Value = 0.86699999999999999
$qol.Range("F9").Value = 0.86099999999999999

# age_covid sheet: update columns D,E,F,G,H (age-specific COVID rates)
$age = $wb.Worksheets.Item("age_covid")

$age.Range("D2").Value = 0.000060475336142076724
$age.Range("E2").Value = 0.00019603300482954038
$age.Range("F2").Value = 0.00011305822354073918
$age.Range("H2").Value = 0.00093023255813953494

$age.Range("D3").Value = 0.0002419013445683069
$age.Range("E3").Value = 0.0006950261080320068
$age.Range("F3").Value = 0.000000012782162073571416
$age.Range("H3").Value = 0.0018604651162790699

$age.Range("D4").Value = 0.0014715665127905337
$age.Range("E4").Value = 0.0040632295546486558
$age.Range("F4").Value = 0.0010175240118666525
$age.Range("H4").Value = 0.0027906976744186047

$age.Range("D5").Value = 0.0041526397484226013
$age.Range("E5").Value = 0.012216420437331812
$age.Range("F5").Value = 0.0016958733531110878
$age.Range("H5").Value = 0.0037209302325581397

$age.Range("D6").Value = 0.013909327312677647
$age.Range("E6").Value = 0.03323650490973571
$age.Range("F6").Value = 0.0053137365064147412
$age.Range("G6").Value = 0.01984126984126984
$age.Range("H6").Value = 0.013023255813953489

$age.Range("D7").Value = 0.045719354123410001
$age.Range("E7").Value = 0.084441216830324517
$age.Range("F7").Value = 0.02340305227293301
$age.Range("G7").Value = 0.023809523809523808
$age.Range("H7").Value = 0.043720930232558138

$age.Range("D8").Value = 0.096679904045799986
$age.Range("E8").Value = 0.16394596617539609
$age.Range("F8").Value = 0.071452797277747157
$age.Range("G8").Value = 0.083333333333333329
$age.Range("H8").Value = 0.10232558139534884

$age.Range("D9").Value = 0.22504888423004818
$age.Range("E9").Value = 0.23663411330707679
$age.Range("F9").Value = 0.18111927411226417
$age.Range("G9").Value = 0.23412698412698413
$age.Range("H9").Value = 0.23069767441860464

$age.Range("D10").Value = 0.39520632168847136
$age.Range("E10").Value = 0.29862509578885466
$age.Range("F10").Value = 0.35794233572998024
$age.Range("G10").Value = 0.34920634920634919
$age.Range("H10").Value = 0.35627906976744184

$age.Range("D11").Value = 0.21750962565766929
$age.Range("E11").Value = 0.16594639388377025
$age.Range("F11").Value = 0.35794233572998024
$age.Range("G11").Value = 0.28968253968253971
$age.Range("H11").Value = 0.24465116279069768

# ---------------------------------------------------------------------------
# Sheet view / selection bookkeeping, matching the author's final state:
#   - female_LT:  scroll/selection moves to B15 (drops the old topLeftCell)
#   - age_covid:  selection moves to I6, no longer the active tab
#   - qol_norm:   selection moves to F12, becomes the active tab
# Order matters: the last sheet Activate()'d becomes the workbook's active
# tab, so qol_norm is activated last.
# ---------------------------------------------------------------------------
$female = $wb.Worksheets.Item("female_LT")
[void]$female.Activate()
[void]$female.Range("B15").Select()

[void]$age.Activate()
[void]$age.Range("I6").Select()

[void]$qol.Activate()
[void]$qol.Range("F12").Select()

Write-Host "Edit complete"
